$d = $word.ActiveDocument

# 1) Heading1 paragraph spacing-before: 240 -> 183 twips (12pt -> 9.15pt)
$d.Paragraphs(1).Format.SpaceBefore = 9.15

# 2) Table cell margin left: 54 -> 53 twips (2.7pt -> 2.65pt)
$tbl = $d.Tables(1)
$tbl.LeftPadding = 2.65

# 3) Merge split runs back into single runs via Find/Replace (same text),
#    walking forward through the document so duplicate phrases are matched
#    in document order, once each.
$r = $d.Content
$r.Collapse(1)  # wdCollapseStart
$ok = $r.Find.Execute("Stwórz puste repozytorium w aktualnym katalogu.", $true, $false, $false, $false, $false, $true, 0, $false, "Stwórz puste repozytorium w aktualnym katalogu.", 2)
if (-not $ok) { Write-Host "FAILED match #0: Stwórz puste repozytorium w aktualnym katalogu." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Dodaje wszystkie zmiany ze ", $true, $false, $false, $false, $false, $true, 0, $false, "Dodaje wszystkie zmiany ze ", 2)
if (-not $ok) { Write-Host "FAILED match #1: Dodaje wszystkie zmiany ze " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Tworzy commita ze zmianami aktualnie w ", $true, $false, $false, $false, $false, $true, 0, $false, "Tworzy commita ze zmianami aktualnie w ", 2)
if (-not $ok) { Write-Host "FAILED match #2: Tworzy commita ze zmianami aktualnie w " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Pokazuje, które pliki się zmieniły od ostatniego commita, a które zostały dodane lub usunięte i czy są lokalnie jakieś commity do wypushowania lub pobrania.", $true, $false, $false, $false, $false, $true, 0, $false, "Pokazuje, które pliki się zmieniły od ostatniego commita, a które zostały dodane lub usunięte i czy są lokalnie jakieś commity do wypushowania lub pobrania.", 2)
if (-not $ok) { Write-Host "FAILED match #3: Pokazuje, które pliki się zmieniły od ostatniego commita, a które zostały dodane lub usunięte i czy są lokalnie jakieś commity do wypushowania lub pobrania." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Tworzy nowy branch o nazwie ", $true, $false, $false, $false, $false, $true, 0, $false, "Tworzy nowy branch o nazwie ", 2)
if (-not $ok) { Write-Host "FAILED match #4: Tworzy nowy branch o nazwie " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Przełącza się", $true, $false, $false, $false, $false, $true, 0, $false, "Przełącza się", 2)
if (-not $ok) { Write-Host "FAILED match #5: Przełącza się" }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("git checkout -b <name>", $true, $false, $false, $false, $false, $true, 0, $false, "git checkout -b <name>", 2)
if (-not $ok) { Write-Host "FAILED match #6: git checkout -b <name>" }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Resetuje", $true, $false, $false, $false, $false, $true, 0, $false, "Resetuje", 2)
if (-not $ok) { Write-Host "FAILED match #7: Resetuje" }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute(" (czasami bezpowrotnie - zależnie od scenariusza).", $true, $false, $false, $false, $false, $true, 0, $false, " (czasami bezpowrotnie - zależnie od scenariusza).", 2)
if (-not $ok) { Write-Host "FAILED match #8:  (czasami bezpowrotnie - zależnie od scenariusza)." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Resetuje", $true, $false, $false, $false, $false, $true, 0, $false, "Resetuje", 2)
if (-not $ok) { Write-Host "FAILED match #9: Resetuje" }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute(" (jest to domyślny sposób resetowania), nie wymaga przełącznika.", $true, $false, $false, $false, $false, $true, 0, $false, " (jest to domyślny sposób resetowania), nie wymaga przełącznika.", 2)
if (-not $ok) { Write-Host "FAILED match #10:  (jest to domyślny sposób resetowania), nie wymaga przełącznika." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Resetuje repozytorium do wersji ", $true, $false, $false, $false, $false, $true, 0, $false, "Resetuje repozytorium do wersji ", 2)
if (-not $ok) { Write-Host "FAILED match #11: Resetuje repozytorium do wersji " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Usuwa zmiany naniesione na ", $true, $false, $false, $false, $false, $true, 0, $false, "Usuwa zmiany naniesione na ", 2)
if (-not $ok) { Write-Host "FAILED match #12: Usuwa zmiany naniesione na " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Usuwa ", $true, $false, $false, $false, $false, $true, 0, $false, "Usuwa ", 2)
if (-not $ok) { Write-Host "FAILED match #13: Usuwa " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("<file>. A", $true, $false, $false, $false, $false, $true, 0, $false, "<file>. A", 2)
if (-not $ok) { Write-Host "FAILED match #14: <file>. A" }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute(". Zazwyczaj: ", $true, $false, $false, $false, $false, $true, 0, $false, ". Zazwyczaj: ", 2)
if (-not $ok) { Write-Host "FAILED match #15: . Zazwyczaj: " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Szuka commita, który był wspólnym przodkiem zarówno aktualnego brancha jak i ", $true, $false, $false, $false, $false, $true, 0, $false, "Szuka commita, który był wspólnym przodkiem zarówno aktualnego brancha jak i ", 2)
if (-not $ok) { Write-Host "FAILED match #16: Szuka commita, który był wspólnym przodkiem zarówno aktualnego brancha jak i " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("tworzy commita łączącego ", $true, $false, $false, $false, $false, $true, 0, $false, "tworzy commita łączącego ", 2)
if (-not $ok) { Write-Host "FAILED match #17: tworzy commita łączącego " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute(" nie ma niezależnej historii, następuje połączenie ", $true, $false, $false, $false, $false, $true, 0, $false, " nie ma niezależnej historii, następuje połączenie ", 2)
if (-not $ok) { Write-Host "FAILED match #18:  nie ma niezależnej historii, następuje połączenie " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute(", czyli dodanie commitów, bez utworzenia commita mergującego.", $true, $false, $false, $false, $false, $true, 0, $false, ", czyli dodanie commitów, bez utworzenia commita mergującego.", 2)
if (-not $ok) { Write-Host "FAILED match #19: , czyli dodanie commitów, bez utworzenia commita mergującego." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Pozwala zmieniać kolejność commitów, łączyć je  ze sobą, usuwać, modyfikować, zaczynając od ", $true, $false, $false, $false, $false, $true, 0, $false, "Pozwala zmieniać kolejność commitów, łączyć je  ze sobą, usuwać, modyfikować, zaczynając od ", 2)
if (-not $ok) { Write-Host "FAILED match #20: Pozwala zmieniać kolejność commitów, łączyć je  ze sobą, usuwać, modyfikować, zaczynając od " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Tworzy nowego commita, który odwraca zmiany z ", $true, $false, $false, $false, $false, $true, 0, $false, "Tworzy nowego commita, który odwraca zmiany z ", 2)
if (-not $ok) { Write-Host "FAILED match #21: Tworzy nowego commita, który odwraca zmiany z " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Wysyła lokalne commity do zdalnego repozytorium.", $true, $false, $false, $false, $false, $true, 0, $false, "Wysyła lokalne commity do zdalnego repozytorium.", 2)
if (-not $ok) { Write-Host "FAILED match #22: Wysyła lokalne commity do zdalnego repozytorium." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Pobiera commity do lokalnego repozytorium", $true, $false, $false, $false, $false, $true, 0, $false, "Pobiera commity do lokalnego repozytorium", 2)
if (-not $ok) { Write-Host "FAILED match #23: Pobiera commity do lokalnego repozytorium" }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Wykonuje ", $true, $false, $false, $false, $false, $true, 0, $false, "Wykonuje ", 2)
if (-not $ok) { Write-Host "FAILED match #24: Wykonuje " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("git merge.", $true, $false, $false, $false, $false, $true, 0, $false, "git merge.", 2)
if (-not $ok) { Write-Host "FAILED match #25: git merge." }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute("Wykonuje ", $true, $false, $false, $false, $false, $true, 0, $false, "Wykonuje ", 2)
if (-not $ok) { Write-Host "FAILED match #26: Wykonuje " }
$r.Collapse(0)  # wdCollapseEnd
$ok = $r.Find.Execute(", następnie rebasuje lokalne niewypushowane commity, nakładając je na wierzch commitów pobranych ze zdalnego repozytorium.", $true, $false, $false, $false, $false, $true, 0, $false, ", następnie rebasuje lokalne niewypushowane commity, nakładając je na wierzch commitów pobranych ze zdalnego repozytorium.", 2)
if (-not $ok) { Write-Host "FAILED match #27: , następnie rebasuje lokalne niewypushowane commity, nakładając je na wierzch commitów pobranych ze zdalnego repozytorium." }
$r.Collapse(0)  # wdCollapseEnd

# 4) Add a top border (single, 2, #000000) to the tcBorders of the 2nd..6th row cells
for ($i = 2; $i -le $tbl.Rows.Count; $i++) {
  $row = $tbl.Rows($i)
  $rb = $row.Borders
  $top = $rb.Item(-1)   # wdBorderTop
  $top.LineStyle = 1    # wdLineStyleSingle
  $top.LineWidth = 1    # wdLineWidth025pt -> sz=2
  $top.Color = 0        # black
  $rb.DistanceFromTop = 0
}

Write-Host "done"